# Update countries & provincias Spain
# - Reorders 3 country rows (Ucrania/Bolivia, Senegal/Malasia, Croacia/Hungria)
#   so the country name that got *new* data lands on the earlier of the two
#   adjacent rows, while the other country keeps the previous occupant's old
#   numbers on the later row.
# - Refreshes the case/death counters for a number of other countries.
# - Bumps the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6 (India) - counters refresh only
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 2).Value = 1080194
$ws.Cells.Item(6, 3).Value = 2330
$ws.Cells.Item(6, 4).Value = 678709
$ws.Cells.Item(6, 5).Value = 374643
$ws.Cells.Item(6, 7).Value = 14
$ws.Cells.Item(6, 8).Value = 26842

# ---------------------------------------------------------------------------
# Row 25 (Catar) - counters refresh only
# ---------------------------------------------------------------------------
$ws.Cells.Item(25, 2).Value = 106648
$ws.Cells.Item(25, 3).Value = 340
$ws.Cells.Item(25, 4).Value = 103377
$ws.Cells.Item(25, 5).Value = 3114
$ws.Cells.Item(25, 7).Value = 3
$ws.Cells.Item(25, 8).Value = 157

# ---------------------------------------------------------------------------
# Row 35 (Bielorrusia) - counters refresh only
# ---------------------------------------------------------------------------
$ws.Cells.Item(35, 2).Value = 66095
$ws.Cells.Item(35, 3).Value = 142
$ws.Cells.Item(35, 4).Value = 58204
$ws.Cells.Item(35, 5).Value = 7392
$ws.Cells.Item(35, 7).Value = 4
$ws.Cells.Item(35, 8).Value = 499

# ---------------------------------------------------------------------------
# Row 37 (Kuwait) - counters refresh only
# ---------------------------------------------------------------------------
$ws.Cells.Item(37, 2).Value = 59204
$ws.Cells.Item(37, 3).Value = 300
$ws.Cells.Item(37, 4).Value = 49687
$ws.Cells.Item(37, 5).Value = 9109
$ws.Cells.Item(37, 7).Value = 1
$ws.Cells.Item(37, 8).Value = 408

# ---------------------------------------------------------------------------
# Rows 38-39: Ucrania moves up to row 38 with fresh counters; Bolivia moves
# down to row 39 keeping the counters row 38 used to have.
# ---------------------------------------------------------------------------
$ws.Cells.Item(38, 1).Value = "Ucrania"
$ws.Cells.Item(38, 2).Value = 58842
$ws.Cells.Item(38, 3).Value = 731
$ws.Cells.Item(38, 4).Value = 30879
$ws.Cells.Item(38, 5).Value = 26478
$ws.Cells.Item(38, 7).Value = 8
$ws.Cells.Item(38, 8).Value = 1485

$ws.Cells.Item(39, 1).Value = "Bolivia"
$ws.Cells.Item(39, 2).Value = 58138
$ws.Cells.Item(39, 3).Value = 2036
$ws.Cells.Item(39, 4).Value = 18200
$ws.Cells.Item(39, 5).Value = 37832
$ws.Cells.Item(39, 7).Value = 57
$ws.Cells.Item(39, 8).Value = 2106

# ---------------------------------------------------------------------------
# Rows 82-83: Senegal moves up to row 82 with fresh counters; Malasia moves
# down to row 83 keeping the counters row 82 used to have.
# ---------------------------------------------------------------------------
$ws.Cells.Item(82, 1).Value = "Senegal"
$ws.Cells.Item(82, 2).Value = 8810
$ws.Cells.Item(82, 3).Value = 141
$ws.Cells.Item(82, 4).Value = 5948
$ws.Cells.Item(82, 5).Value = 2695
$ws.Cells.Item(82, 7).Value = 4
$ws.Cells.Item(82, 8).Value = 167

$ws.Cells.Item(83, 1).Value = "Malasia"
$ws.Cells.Item(83, 2).Value = 8779
$ws.Cells.Item(83, 3).Value = 15
$ws.Cells.Item(83, 4).Value = 8553
$ws.Cells.Item(83, 5).Value = 103
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 123

# ---------------------------------------------------------------------------
# Row 85 (Estado de Palestina) - partial counters refresh only
# ---------------------------------------------------------------------------
$ws.Cells.Item(85, 5).Value = 6568
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 60

# ---------------------------------------------------------------------------
# Row 87 (Bosnia y Herzegovina) - counters refresh only
# ---------------------------------------------------------------------------
$ws.Cells.Item(87, 2).Value = 8340
$ws.Cells.Item(87, 3).Value = 179
$ws.Cells.Item(87, 4).Value = 3669
$ws.Cells.Item(87, 5).Value = 4422
$ws.Cells.Item(87, 7).Value = 3
$ws.Cells.Item(87, 8).Value = 249

# ---------------------------------------------------------------------------
# Row 90 (Madagascar) - counters refresh only
# ---------------------------------------------------------------------------
$ws.Cells.Item(90, 2).Value = 7049
$ws.Cells.Item(90, 3).Value = 200
$ws.Cells.Item(90, 4).Value = 3498
$ws.Cells.Item(90, 5).Value = 3492
$ws.Cells.Item(90, 7).Value = 4
$ws.Cells.Item(90, 8).Value = 59

# ---------------------------------------------------------------------------
# Rows 99-100: Croacia moves up to row 99 with fresh counters; Hungria moves
# down to row 100 keeping the counters row 99 used to have.
# ---------------------------------------------------------------------------
$ws.Cells.Item(99, 1).Value = "Croacia"
$ws.Cells.Item(99, 2).Value = 4345
$ws.Cells.Item(99, 3).Value = 92
$ws.Cells.Item(99, 4).Value = 3018
$ws.Cells.Item(99, 5).Value = 1207
$ws.Cells.Item(99, 8).Value = 120

$ws.Cells.Item(100, 1).Value = "Hungria"
$ws.Cells.Item(100, 2).Value = 4333
$ws.Cells.Item(100, 3).Value = 18
$ws.Cells.Item(100, 4).Value = 3223
$ws.Cells.Item(100, 5).Value = 514
$ws.Cells.Item(100, 8).Value = 596

# ---------------------------------------------------------------------------
# Row 143 (Burkina Faso) - counters refresh only
# ---------------------------------------------------------------------------
$ws.Cells.Item(143, 2).Value = 1052
$ws.Cells.Item(143, 3).Value = 5
$ws.Cells.Item(143, 4).Value = 901
$ws.Cells.Item(143, 5).Value = 98

# ---------------------------------------------------------------------------
# Bump the "last updated" timestamp shown in A1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 14:24"
